$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: drop the wrap (keep the horizontal centering) ---
$ws.Range("A1:H1").WrapText = $false

# --- New lookup/status columns (I:L) ---
$ws.Range("I1").Value = "Name"
$ws.Range("J1").Value = "Price"
$ws.Range("K1").Value = "Availability"
$ws.Range("L1").Value = "Purchased"

# --- Normalize the new data block's vertical alignment back to default ---
$ws.Range("A2:H6").VerticalAlignment = -4108
$ws.Range("A2:H6").VerticalAlignment = -4107

# --- Row 2: Irwin Vise-Grip pliers (listing still undefined) ---
$ws.Range("C2").VerticalAlignment = -4108
$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.walmart.com/ip/IRWIN-VISE-GRIP-2076709-Fast-Release-Locking-Plier-Set-6Pc/40751123")
$ws.Range("I2").Value = "Undefined"
$ws.Range("J2").Value = "Undefined"
$ws.Range("K2").Value = "Something wrong with this listing. Check it by yourself!"
$ws.Range("L2").Value = "Undefined"

# --- Row 3: Chamberlain wireless motion alert (listing still undefined) ---
$ws.Range("C3").VerticalAlignment = -4108
$ws.Hyperlinks.Add($ws.Range("C3"), "https://www.walmart.com/ip/CHAMBERLAIN-CWA2000-Wireless-Motion-Alert-System/14554538")
$ws.Range("I3").Value = "Undefined"
$ws.Range("J3").Value = "Undefined"
$ws.Range("K3").Value = "Something wrong with this listing. Check it by yourself!"
$ws.Range("L3").Value = "Undefined"

# --- Row 4: Air Purifier (resolved, different vertical alignment than the rest) ---
$ws.Hyperlinks.Add($ws.Range("C4"), "https://www.walmart.com/ip/Air-Purifier-Hepa-Carbon-Ionic-Ozone-Generator-Cleaner-UV-C-with-Remote/179032464")
$ws.Range("I4").Value = "Air Purifier Hepa Carbon Ionic Ozone Generator Cleaner UV-C, with Remote"
$ws.Range("J4").Value = "97,93"
$ws.Range("K4").Value = "Availible"
$ws.Range("L4").Value = "Undefined"

# --- Row 5: 4PC self-adjusting pipe wrench (resolved) ---
$ws.Range("C5").VerticalAlignment = -4108
$ws.Hyperlinks.Add($ws.Range("C5"), "https://www.walmart.com/ip/4PC-Self-Adjusting-Quick-Release-Pipe-Wrench-Drop-Forge-Plumbing/193316668")
$ws.Range("I5").Value = "4PC Self-Adjusting Quick Release Pipe Wrench Drop Forge Plumbing"
$ws.Range("J5").Value = "49,93"
$ws.Range("K5").Value = "Availible"
$ws.Range("L5").Value = "Undefined"

# --- Row 6: Torque multiplier wrench (resolved) ---
$ws.Range("C6").VerticalAlignment = -4108
$ws.Hyperlinks.Add($ws.Range("C6"), "https://www.walmart.com/ip/Torque-Multiplier-Wrench-Lug-Nut-Remover-with-4-cr-v-sockets/102917115")
$ws.Range("I6").Value = "Torque Multiplier Wrench Lug Nut Remover, with 4 cr-v sockets"
$ws.Range("J6").Value = "53,93"
$ws.Range("K6").Value = "Availible"
$ws.Range("L6").Value = "Undefined"

# --- Restore the cursor where the user last clicked before saving ---
$ws.Range("C12").Select() | Out-Null
